$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 298, pushing existing rows 298:304 down to 299:305
$ws.Rows("298:298").Insert()

# Populate the newly inserted row 298 with the new record
$ws.Range("A298").Value = 10
$ws.Range("B298").Value = "Vega Modelo de Temuco"
$ws.Range("C298").Value = "La Araucanía"
$ws.Range("D298").Value = 44890
$ws.Range("E298").Value = 9
$ws.Range("F298").Value = "Fruta"
$ws.Range("G298").Value = 100103
$ws.Range("H298").Value = "Frutos de hueso (carozo)"
$ws.Range("I298").Value = 100103004
$ws.Range("J298").Value = "Durazno"
$ws.Range("K298").Value = "Early Majestic"
$ws.Range("L298").Value = "Primera"
$ws.Range("M298").Value = 100
$ws.Range("N298").Value = 15000
$ws.Range("O298").Value = 15000
$ws.Range("P298").Value = 15000
$ws.Range("Q298").Value = '$/bandeja 10 kilos granel'
$ws.Range("R298").Value = "Provincia de Limarí"
$ws.Range("S298").Value = 1500
$ws.Range("T298").Value = 10
